$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "hours worked" header in D1, matching the existing header style (A1:C1) ---
$ws.Range("D1").Value = "hours worked"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# --- "hours worked" values for the existing rows (2-17) ---
$hours = @{
    2 = 3; 3 = 3; 4 = 3; 5 = 3; 6 = 3; 7 = 3; 8 = 3; 9 = 3; 10 = 3
    11 = 2; 12 = 2
    13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1
}
foreach ($r in $hours.Keys) {
    $ws.Cells.Item($r, 4).Value = $hours[$r]
}

# --- Three new weekly rows (18-20), formatted like row 17 (Date / Time in / Time out) ---
$ws.Range("A17:C17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)

$ws.Range("A18").Value = 43430
$ws.Range("B18").Value = 0.375
$ws.Range("C18").Value = 0.458333333333333
$ws.Range("D18").Value = 2

$ws.Range("A19").Value = 43437
$ws.Range("B19").Value = 0.375
$ws.Range("C19").Value = 0.458333333333333
$ws.Range("D19").Value = 2

$ws.Range("A20").Value = 43444
$ws.Range("B20").Value = 0.375
$ws.Range("C20").Value = 0.458333333333333
$ws.Range("D20").Value = 2

# --- Totals row 21: bold "hours worked" label + the summed total ---
$ws.Range("D21").Value = "hours worked"
$ws.Range("E21").Value = 42
$ws.Range("D21:E21").Font.Bold = $true

# --- Widen column D so the new header/label text fits ---
$ws.Columns.Item(4).ColumnWidth = 13.65

# --- Restore the selection left behind by the author's last edit ---
$ws.Range("G16").Select() | Out-Null
